# Applies the "Quantum Entanglement" -> "The Profound Significance of Arts in
# Society" content swap described by the commit diff.
#
# Strategy: every change in the diff is, at the text level, either
#   (a) a straight 1:1 run text replacement, or
#   (b) several adjacent runs collapsing into a single run (the diff removes
#       the extra <w:r> elements entirely and keeps only the first run's
#       rPr/formatting).
# Using Find.Execute with the concatenation of the old runs' text as the
# search string and the new sentence as the replacement reproduces both
# cases faithfully: Word collapses the matched range (regardless of how many
# runs/paragraph-internal runs it spans) into a single run that carries the
# formatting of the range's first run.

$d = $word.ActiveDocument

function Replace-Text {
    param(
        [string]$Old,
        [string]$New,
        [bool]$WholeWord = $false
    )
    $d.Content.Find.Execute($Old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $New, 2, $false, $false, $false, `
                             $false) | Out-Null
}

# --- Title ---------------------------------------------------------------
$d.Content.Find.Execute("Quantum Entanglement: Unveiling the Mysteries of Interconnectedness", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "The Profound Significance of Arts in Society", 2) | Out-Null

# --- Byline: 5 runs -> 1 run ----------------------------------------------
$d.Content.Find.Execute("Dr. Richard E. Jester", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Cheryl Manning", 2) | Out-Null

# --- Email address ---------------------------------------------------------
$d.Content.Find.Execute("rejester@quanticoinstitute", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "cherylmnnng@icloud", 2) | Out-Null

# "com" is a short, common substring (computing/computers elsewhere), so we
# must match it as a whole word to hit only the isolated "com" run.
$d.Content.Find.Execute("com", `
    $true, $true, $false, $false, $false, $true, 1, $false, `
    "net", 2) | Out-Null

# --- Body paragraph (first <w:br/> block) ----------------------------------
$d.Content.Find.Execute("In the realm of quantum mechanics, there exists a phenomenon that defies classical understanding and challenges our perception of reality: quantum entanglement", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "In the grand tapestry of human civilization, the arts find their home as a vibrant expression of our collective consciousness", 2) | Out-Null

$d.Content.Find.Execute(" This extraordinary phenomenon unveils a universe where particles, separated by vast distances, exhibit a profound interconnectedness, defying the constraints of space and time", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " From the echoing chambers of music to the canvas of painted stories, arts have forever been the vessel through which we navigate the complexities of life and make sense of our place within it", 2) | Out-Null

$d.Content.Find.Execute(" Imagine a universe where the destiny of one particle is inextricably linked to that of another, regardless of the distance that separates them. Such is the enigmatic realm of quantum entanglement, a realm where the boundaries of reality blur and the conventional notions of locality and causality are rendered obsolete. This profound interconnectedness has sparked a plethora of investigations, unraveling the mysteries of quantum entanglement and its implications for our understanding of the universe", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " As we delve into the world of arts, a symphony of colors, sounds, and narratives unfolds, painting pictures of the human spirit and enriching our understanding of the cosmos we call home", 2) | Out-Null

# --- Body paragraph (second <w:br/> block) ---------------------------------
$d.Content.Find.Execute("The exploration of quantum entanglement has ignited a scientific revolution, propelling us into a new era of discovery and innovation", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Arts, in its kaleidoscope of forms, provide a mirror to society, reflecting its trials, tribulations, triumphs, and aspirations", 2) | Out-Null

$d.Content.Find.Execute(" This enigmatic phenomenon has the potential to reshape various fields, ranging from cryptography to computing", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " It becomes a pulpit for voices unheard, a solace for hearts burdened by sorrow, and a window to the soul of humanity", 2) | Out-Null

$d.Content.Find.Execute(" The prospect of harnessing the power of quantum entanglement to construct unbreakable encryption codes has the potential to revolutionize the realm of information security. Moreover, quantum entanglement holds the promise of enabling the development of ultra-fast quantum computers, capable of performing calculations that are intractable for classical computers, heralding a new era of computational prowess. As we delve deeper into the intricacies of quantum entanglement, we unlock the mysteries of interconnectedness, expanding our comprehension of the universe and paving the way for transformative advancements in technology and scientific understanding", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " From the timeless sculptures that narrate tales of ancient civilizations to the evocative lyrics that capture the yearnings of our hearts, arts hold the power to connect us across time and space, forging an enduring bond between people of diverse backgrounds and experiences", 2) | Out-Null

# --- Body paragraph (third <w:br/> block) -----------------------------------
$d.Content.Find.Execute("The profound implications of quantum entanglement extend far beyond the realm of physics, reaching into the depths of philosophy, blurring the boundaries between science and spirituality", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "As vessels of cultural memory, arts bestow upon us a precious understanding of our roots, facilitating a dialogue between generations", 2) | Out-Null

$d.Content.Find.Execute(" This enigmatic phenomenon challenges our fundamental assumptions about reality, inviting us to contemplate the interconnectedness of all things", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " Through stories, dances, and melodic tunes, arts weave the fabric of our collective legacy, preserving customs, traditions, and values that define us as a people", 2) | Out-Null

$d.Content.Find.Execute(" It raises profound questions about the nature of consciousness, the relationship between mind and matter, and the interconnectedness of the universe. The exploration of quantum entanglement propels us into a realm where the boundaries of science and philosophy intersect, encouraging us to question the fundamental nature of reality and our place within it. This journey of discovery not only enriches our scientific understanding but also invites us to ponder the deepest mysteries of existence, beckoning us to seek a deeper connection with the universe and all that it holds", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " In the archival chambers of arts, we find clues to our origins, learn from the wisdom of our ancestors, and pass on lessons to those who come after us, ensuring the flame of our heritage never flickers out", 2) | Out-Null

# --- Summary paragraph -------------------------------------------------------
$d.Content.Find.Execute("Quantum entanglement, a captivating phenomenon in the realm of quantum mechanics, unveils a universe where distant particles exhibit profound interconnectedness, transcending the constraints of space and time", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "In essence, arts are a mirror to society, reflecting our trials, our triumphs, and our collective consciousness", 2) | Out-Null

$d.Content.Find.Execute(" This enigmatic phenomenon has ignited a scientific revolution, with implications that reverberate across diverse fields, from cryptography to computing", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " They serve as a window to the soul of humanity, connecting us across time and space, forging enduring bonds between diverse cultures", 2) | Out-Null

$d.Content.Find.Execute(" The potential to harness the power of quantum entanglement for unbreakable encryption codes and ultra-fast quantum computers holds immense promise for technological advancements", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " As vessels of cultural memory, arts preserve our customs, traditions, and values, facilitating a dialogue between generations", 2) | Out-Null

$d.Content.Find.Execute(" Beyond its practical applications, quantum entanglement challenges our fundamental assumptions about reality, inviting us to contemplate the interconnectedness of all things and ponder the deepest mysteries of existence. As we delve deeper into the intricacies of quantum entanglement, we embark on a journey that not only expands our scientific understanding but also enriches our philosophical ponderings, propelling us towards a more profound comprehension of the universe and our place within it", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " Through arts, we learn from our ancestors and pass on lessons to those who come after us, ensuring the flame of heritage never flickers out", 2) | Out-Null

# --- Trailing empty paragraph added after the Summary paragraph -------------
$d.Paragraphs.Add() | Out-Null
